$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @('File', 'Review date'),
    @('Routine_Care/Nursing for Arterial and Central Venous Lines.pdf', '2011-03-12'),
    @('Routine_Care/VTE_Prevention/TED Stocking Sizing.pdf', '2011-06-12'),
    @('Routine_Care/Faecal  incontinence skin care.pdf', '2011-07-12'),
    @('Breathing(Respiratory)/Equipment/IPPB using an ICU Ventilator.pdf', '2011-12-12'),
    @('Drugs/heparin_critical_care_only.pdf', '2014-09-12'),
    @('Routine_Care/Invasive Flush Systems.pdf', '2014-10-12'),
    @('GI_Liver_and_Transplant/Pancreatic Irrigation.pdf', '2014-11-12'),
    @('Neurological/Thiopentone levels.pdf', '2014-12-12'),
    @('Infection_and_sepsis/Ebola/Ebola.pdf', '2015-03-12'),
    @('Breathing(Respiratory)/Equipment/AMBU AScope.pdf', '2016-04-12'),
    @('Cardiovascular/EZ-IO Intraosseus Access Device_pub_em.pdf', '2016-08-12'),
    @('Routine_Care/Central venous catheter removal.pdf', '2017-01-12'),
    @('Routine_Care/Tracheostomy_nursing_care.pdf', '2017-08-12'),
    @('ECLS/Extra Corporeal Carbon Dioxide Removal.pdf', '2017-10-12'),
    @('Airway/Critical care extubation checklist.pdf', '2017-10-12'),
    @('Airway/Tracheostomy_Laryngectomy/Hospital_in-patients_with_a_Tracheostomy.pdf', '2017-11-12'),
    @('Drugs/diazepam_diazemuls.pdf', '2018-07-12'),
    @('End_of_life_care/Reasons to report a death to PF.pdf', '2019-05-12'),
    @('Breathing(Respiratory)/salbutamol and ipratroprium MDI.pdf', '2019-05-12'),
    @('Neurological/SOP -  Femoral site care.pdf', '2019-06-12'),
    @('Airway/Tracheostomy_Laryngectomy/Tracheostomy change in Critical Care.pdf', '2019-06-12'),
    @('Airway/Tracheostomy_Laryngectomy/Tracheostomy suctioning cleaning guideline.pdf', '2019-06-12'),
    @('Drugs/sodium_bicarbonate.pdf', '2019-08-12'),
    @('Post_op_care/Anticoagulation antiplatelet agents and epidural analgesia.pdf', '2019-09-12'),
    @('Post_op_care/Epidural top-up.pdf', '2020-01-12'),
    @('Covid-19/SJH/SJH COVID19 ED Intubation Action Card.pdf', '2020-03-12'),
    @('Diabetes_and_Glucose/Hyperosmolar Hyperglycaemic State.pdf', '2020-03-12'),
    @('Covid-19/SJH/SJH COVID19 ITU Intubation Action Card.pdf', '2020-03-12'),
    @('Covid-19/WGH/CoVid intubation checklist WGH.pdf', '2020-03-12'),
    @('Drugs/heparin for Haemofiltration.pdf', '2020-03-12'),
    @('Airway/Emergency intubation checklist_em_pub.pdf', '2020-03-12'),
    @('Drugs/fentanyl.pdf', '2020-04-12'),
    @('Airway/Tracheostomy_Laryngectomy/Tracheostomy guideline.pdf', '2020-05-12'),
    @('Covid-19/WGH/WGH_CT_Transfer_May.pdf', '2020-07-12'),
    @('Cardiovascular/GJNH Acute Heart Failure Referral Form.pdf', '2020-08-12'),
    @('Organ_donation/Donation after circulatory death.pdf', '2020-11-12'),
    @('Airway/Percutaneous tracheostomy checklist.pdf', '2021-02-12'),
    @('Delirium/Risk assessment posi mit.pdf', '2021-05-12'),
    @('Infection_and_sepsis/SOP Ultrasound Cleaning.pdf', '2021-05-12'),
    @('Airway/Tracheostomy_Laryngectomy/Tracheostomy safety box contents.pdf', '2021-06-12'),
    @('Airway/McGrath Mac.pdf', '2021-06-12'),
    @('Airway/Anticipated difficult airway tool.pdf', '2021-06-12'),
    @('GI_Liver_and_Transplant/Abdominal pressure measurement.pdf', '2021-06-12'),
    @('GI_Liver_and_Transplant/Treatment of constipation.pdf', '2021-06-12'),
    @('Neurological/Sub arachnoid haemorrhage management.pdf', '2021-06-12'),
    @('Breathing(Respiratory)/HFNO.pdf', '2021-06-12'),
    @('End_of_life_care/Documentation following death.pdf', '2021-09-12'),
    @('Drugs/zanamivir.pdf', '2021-12-12'),
    @('Routine_Care/bBraun Spaceplus Failure EMERGENCY ACTION CARD_em.pdf', '2022-01-12'),
    @('Drugs/insulin.pdf', '2022-03-12'),
    @('Breathing(Respiratory)/Equipment/HFNO Set Up.pdf', '2022-03-12'),
    @('Breathing(Respiratory)/Inhaled Nitrous Oxide.pdf', '2022-04-12'),
    @('Breathing(Respiratory)/Equipment/APRV.pdf', '2022-05-12'),
    @('Cardiovascular/Steroids for Septic Shock.pdf', '2022-05-12'),
    @('Neurological/SOP for review of Neurosurgical patients in ITU by neurosurgical team.pdf', '2022-06-12'),
    @('Post_op_care/Epidural Haematoma.pdf', '2022-06-12'),
    @('Breathing(Respiratory)/Equipment/T piece Y piece.pdf', '2022-06-12'),
    @('Policies_and_admin/General Critical Care Interaction with HEPMA_pub.pdf', '2022-07-12'),
    @('Drugs/midazolam and thiopental levels.pdf', '2022-08-12'),
    @('Breathing(Respiratory)/Equipment/HFNO through ventilator.pdf', '2022-10-12'),
    @('Post_op_care/Adult Scoliosis Spinal Surgery Post-Op Care.pdf', '2022-11-12'),
    @('Routine_Care/VTE_Prevention/Dalteparin_thromboprophylaxis.pdf', '2022-11-12'),
    @('Post_op_care/Post op care pharyngo-laryngo-oesphagectomy PLOG.pdf', '2022-12-12'),
    @('GI_Liver_and_Transplant/Nasogastric feeding protocol.pdf', '2023-01-12'),
    @('Drugs/Antibiotic doses in CVVHD.pdf', '2023-02-12'),
    @('ECLS/RIE ECLS Anti Xa Protocol.pdf', '2023-04-12'),
    @('GI_Liver_and_Transplant/Jejunostomy feeding protocol.pdf', '2023-04-12'),
    @('GI_Liver_and_Transplant/Nasojejunal feeding protocol.pdf', '2023-04-12'),
    @('Drugs/vasopressin_sepsis.pdf', '2023-05-12'),
    @('Drugs/vasopressin organ donation.pdf', '2023-05-12'),
    @('Infection_and_sepsis/Winter Infections Stepdown Guidance.pdf', '2023-05-12'),
    @('Transfer/ACCP Transfers.pdf', '2023-06-12'),
    @('Covid-19/videos/Donning and Doffing Video.pdf', '2023-06-12'),
    @('Breathing(Respiratory)/Equipment/Ventilators Circuits Filters and Closed Suction - Set up and Maintenance.pdf', '2023-07-12'),
    @('Drugs/piperacillin_tazobactam extended_infusion.pdf', '2023-07-12'),
    @('Infection_and_sepsis/Infection indications for IVIG.pdf', '2023-07-12'),
    @('Procedures/CVC Guidance/Securing CVCs.pdf', '2023-08-12'),
    @('Covid-19/Covid 19 Death Certification Guideline.pdf', '2023-08-12'),
    @('Routine_Care/Video Communication.pdf', '2023-09-12'),
    @('Drugs/isoprenaline.pdf', '2023-10-12'),
    @('Cardiovascular/Cardiogenic Shock.pdf', '2023-10-12'),
    @('Haematology_CAR-T/Haem_ICU_transfer.pdf', '2024-01-12'),
    @('Cardiovascular/Management of hypertension within Critical Care.pdf', '2024-02-12'),
    @('Haematology_CAR-T/Cytokine Release Syndrome.pdf', '2024-03-12'),
    @('Haematology_CAR-T/ICANS.pdf', '2024-03-12'),
    @('Drugs/pancuronium.pdf', '2024-03-12'),
    @('Policies_and_admin/General Critical Care SOP_pub.pdf', '2024-04-12'),
    @('Drugs/milrinone.pdf', '2024-04-12'),
    @('Ethics_and_Law/Care at the End of Life (FICM).pdf', '2024-05-12'),
    @('GI_Liver_and_Transplant/ICU - Upper GI bleeding (Endoscopy guideline).pdf', '2024-05-12'),
    @('Organ_donation/Organ Retrieval SOP.pdf', '2024-05-12'),
    @('End_of_life_care/Palliative extubation & withdrawal of invasive ventilatory support nursing checklist.pdf', '2024-05-12'),
    @('Drugs/clonidine.pdf', '2024-05-12'),
    @('Covid-19/COVID 19 ICM guidance basic goals_June_2022.pdf', '2024-05-12'),
    @('Neurological/Management of traumatic brain injury.pdf', '2024-05-12'),
    @('Infection_and_sepsis/Trip Out of Unit infection guidance.pdf', '2024-05-12'),
    @('End_of_life_care/CMO & NRS Guidance for Doctors completing MCCD.pdf', '2024-05-12'),
    @('Ethics_and_Law/DNACPR policy for Scotland.pdf', '2024-05-12'),
    @('Neurological/Critical Care MRI Procedure_pub.pdf', '2024-05-12'),
    @('Airway/Cook Staged Extubation Set.pdf', '2024-06-02'),
    @('Drugs/noradrenaline (central).pdf', '2024-06-12'),
    @('End_of_life_care/End of Life Care Guideline for Critical Care.pdf', '2024-06-12'),
    @('Post_op_care/Epidural hypotension.pdf', '2024-06-26'),
    @('Drugs/dexmedetomidine.pdf', '2024-07-12'),
    @('GI_Liver_and_Transplant/Confirmation of Nasogastric Tube Position.pdf', '2024-07-12'),
    @('Breathing(Respiratory)/Equipment/Passy Muir Valve.pdf', '2024-07-12'),
    @('Drugs/glyceryl_trinitrate.pdf', '2024-07-12'),
    @('GI_Liver_and_Transplant/Fulminant Liver Failure.pdf', '2024-07-12'),
    @('Infection_and_sepsis/Antifungal guidance in critical care.pdf', '2024-07-25'),
    @('Cardiovascular/Intra Aortic Balloon Pump Guideline_pub.pdf', '2024-08-12'),
    @('Cardiovascular/Intra Aortic Balloon Pump Bedside Checks_pub.pdf', '2024-08-12'),
    @('Transfer/Transfer Outdoors to Garden Guideline.pdf', '2024-08-12'),
    @('Breathing(Respiratory)/ARDS Strategy.pdf', '2024-08-15'),
    @('Drugs/adrenaline.pdf', '2024-10-12'),
    @('Drugs/Dobutamine.pdf', '2024-10-12'),
    @('Drugs/hydralazine.pdf', '2024-10-24'),
    @('Post_op_care/Major OMFS Free Flap.pdf', '2024-11-20'),
    @('Drugs/Alteplase for massive PE.pdf', '2024-11-24'),
    @('Drugs/alfentanil.pdf', '2024-11-24'),
    @('Drugs/magnesium.pdf', '2024-12-24'),
    @('Drugs/neostigmine.pdf', '2025-01-12'),
    @('Drugs/ICU - IV drug infusions.pdf', '2025-01-12'),
    @('Drugs/vancomycin continuous infusion.pdf', '2025-01-12'),
    @('Neurological/Intrathecal policy RIE.pdf', '2025-02-12'),
    @('Drugs/labetalol.pdf', '2025-02-12'),
    @('Infection_and_sepsis/Initial investigation and management in unidentified Infections.pdf', '2025-02-12'),
    @('Cardiovascular/Management of Acute Type B Aortic Dissection Guideline.pdf', '2025-03-12'),
    @('Drugs/nimodipine.pdf', '2025-04-12'),
    @('Drugs/potassium.pdf', '2025-04-12'),
    @('Procedures/CVC Guidance/CVC NHL  April 2023.pdf', '2025-04-12'),
    @('Trauma and Burns/Management of Severe Burns v1.6.pdf', '2025-04-12'),
    @('Drugs/salbutamol.pdf', '2025-04-12'),
    @('Routine_Care/ICU Eye Care Guideline.pdf', '2025-05-12'),
    @('Drugs/phenobarbitone.pdf', '2025-05-12'),
    @('Drugs/amiodarone.pdf', '2025-05-12'),
    @('Procedures/Arterial Line insertion for ACCPs.pdf', '2025-05-12'),
    @('Drugs/nicardipine.pdf', '2025-05-12'),
    @('Drugs/phenylephrine.pdf', '2025-06-12'),
    @('Breathing(Respiratory)/Manual Ventilation and MHI.pdf', '2025-06-12'),
    @('Neurological/Ventriculitis Guideline.pdf', '2025-06-12'),
    @('Drugs/morphine.pdf', '2025-06-12'),
    @('Drugs/Epoprostenol.pdf', '2025-06-12'),
    @('Breathing(Respiratory)/Pulmonary_Embolism_and_DVT/Catheter directed thrombolysis of iliofemoral DVT alteplase_pub.pdf', '2025-07-12'),
    @('Cardiovascular/Cardiac Output Monitoring_pub .pdf', '2025-07-12'),
    @('Drugs/calcium.pdf', '2025-07-12'),
    @('Drugs/atracurium.pdf', '2025-08-12'),
    @('Drugs/dalteparin_thromboprophylaxis.pdf', '2025-08-12'),
    @('Airway/Tracheostomy_Laryngectomy/Decannulation Guidline.pdf', '2025-08-12'),
    @('Drugs/Vancomycin continuous infusion Fluid Restricted.pdf', '2025-08-12'),
    @('Policies_and_admin/Pet Visitation.pdf', '2025-10-12'),
    @('Drugs/valproate.pdf', '2025-10-12'),
    @('Procedures/Inadvertent Catheter Placement Guideline.pdf', '2025-10-12'),
    @('Drugs/ketamine_for_status epilepticus.pdf', '2025-11-12'),
    @('Drugs/stress ulcer prophylaxis.pdf', '2025-11-12'),
    @('Drugs/Phosphate.pdf', '2025-11-12'),
    @('Breathing(Respiratory)/Proning Guideline.pdf', '2026-01-12'),
    @('Drugs/Octreotide.pdf', '2026-01-12'),
    @('GI_Liver_and_Transplant/Prokinetics in ICU.pdf', '2026-01-12'),
    @('Drugs/Thiopentone.pdf', '2026-01-12'),
    @('Neurological/Status epilepticus.pdf', '2026-01-12'),
    @('Procedures/ACCPs acquiring initial CVC competencies.pdf', '2026-03-12'),
    @('Post_op_care/Prevention and treatment of paraplegia after major aortic procedures.pdf', '2026-03-12'),
    @('Procedures/ACCP CVC placement following completion of initial competencies.pdf', '2026-03-12'),
    @('Drugs/Rocuronium.pdf', '2026-05-12'),
    @('Delirium/Violence and Agression.pdf', '2026-05-12'),
    @('Drugs/Ketamine_in_asthma.pdf', '2026-05-12'),
    @('Breathing(Respiratory)/Equipment/O2 through Drager setup.pdf', '2026-05-12'),
    @('Post_op_care/Care of the Transgender Patient.pdf', '2026-07-12'),
    @('Drugs/aminophylline.pdf', '2026-08-12'),
    @('Drugs/IVIG in staph and strep toxic shock syndrome.pdf', '2026-08-12'),
    @('Drugs/Phenytoin.pdf', '2026-09-12'),
    @('GI_Liver_and_Transplant/Plasma exchange in Acute Liver Failure.pdf', '2026-11-12'),
    @('Breathing(Respiratory)/Equipment/NIV through Drager Vent Set up in Critical Care.pdf', '2027-01-12'),
    @('Transfer/Transfer Guidelines.pdf', '2027-02-12'),
    @('End_of_life_care/Guideline following Sudden Cardiac Death where death occurs in ICU.pdf', '2027-03-12'),
    @('Breathing(Respiratory)/Equipment/Setting up the NKV-330.pdf', '2028-02-12'),
    @('Breathing(Respiratory)/NKV 330 NIV guidelines.pdf', '2028-03-12'),
    @('Breathing(Respiratory)/T1RF - CPAP Guidelines_Mar24.pdf', '2028-03-12'),
    @('Infection_and_sepsis/Influenza in Critical Care.pdf', '2028-05-12'),
    @('Policies_and_admin/Anticipated Post op flow surgical patients.pdf', '2028-07-12'),
    @('Policies_and_admin/Discharge Home from Critical Care .pdf', '2028-11-12'),
    @('GI_Liver_and_Transplant/Nasal Bridle use in Critical Care.pdf', '2029-02-12'),
    @('GI_Liver_and_Transplant/Insertion of Nasojejunal Tubes.pdf', '2029-02-12'),
    @('GI_Liver_and_Transplant/Pabrinex.pdf', '2029-06-12'),
    @('Cardiovascular/Focused TOE Guideline for General Critical Care_pub.pdf', '2029-06-12'),
    @('Drugs/Pabrinex.pdf', '2029-06-12'),
    @('Neurological/NPI Guideline.pdf', '2029-08-12'),
    @('Transfer/Repatriaiton Checklist for Critical Care.pdf', '2030-11-12'),
    @('Policies_and_admin/Repatriaiton Checklist for Critical Care.pdf', '2030-11-12')
)

# Force column B to Text format first so date-like strings (e.g. "2011-03-12")
# are stored as literal text rather than being auto-converted to date serials,
# matching the shared-string (t="s") representation used throughout the sheet.
$ws.Range("B2:B190").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $pair = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

# Restore the default (Normal) style on column B so no residual text-format
# style is left applied to the cells, keeping formatting identical to the source.
$ws.Range("B2:B190").Style = "Normal"

Write-Output "done"
